$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the rows that correspond to the "dropped" quarters, working from
# the bottom up so earlier deletions don't shift the row numbers of rows
# still queued for deletion.
$rowsToDelete = @(30, 28, 26, 24, 22, 20, 18, 16, 14, 12, 10, 8, 6, 4, 2)

foreach ($r in $rowsToDelete) {
    $ws.Rows.Item($r).Delete()
}
